$d = $word.ActiveDocument

$d.Content.Find.Execute(
    " _part__option_ _N_ (_UE__PT__N_-_CT__N_)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ", _part__option_ _N_ (_UE__N_-_CT__N_)", 2)

$d.Content.Find.Execute(
    "_First name_ _Last name_ (_email_@_address_._ext_)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "_Forename_ _Surname_ (_local-part_@_domain_)", 2)
